$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# RENOVA&DISEÑA S.A. (row 19) now has a sale of 183.17 in the "240X80 PORCELANATO" column (D)
$ws1.Range("D19").Value = 183.17

# Row 29 holds "x de 27" counters per column - column D gains one more non-zero entry
$ws1.Range("D29").Value = "2 de 27"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same client's "junio" (F) sale updated to 183.17
$ws2.Range("F19").Value = 183.17

# Totals row (29) for junio increases by the same amount
$ws2.Range("F29").Value = 3038.56

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 = "240X80 PORCELANATO" group totals (VENTA / POR CUMPLIR / CUMPLIMIENTO)
$ws3.Range("D3").Value = 274.75
$ws3.Range("E3").Value = 2845.3645
$ws3.Range("F3").Value = 0.08805766583245582

# Row 19 = overall TOTAL row
$ws3.Range("D19").Value = 6294.27
$ws3.Range("E19").Value = 17205.73093005039
$ws3.Range("F19").Value = 0.2678412659954948
